$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between row 55 and row 56 ---
# Row 55 originally: Pyunik Yerevan vs Ararat-Armenia
# Row 56 originally: Van vs Alashkert
# After the edit they trade places (the A:E "index/date" columns stay put).

$row55 = @(
    "Van", 2, "Alashkert", 1,
    6.34, "19/10/2023 00:12", 7.32, "20/10/2023 12:59",
    4.55, "19/10/2023 00:12", 4.51, "20/10/2023 12:59",
    1.38, "19/10/2023 00:12", 1.43, "20/10/2023 12:59",
    "https://www.betexplorer.com/football/armenia/premier-league/van-alashkert/S4mfepEL/"
)

$row56 = @(
    "Pyunik Yerevan", 1, "Ararat-Armenia", 1,
    1.85, "19/10/2023 00:12", 1.99, "20/10/2023 12:59",
    3.41, "19/10/2023 00:12", 3.24, "20/10/2023 12:59",
    3.73, "19/10/2023 00:12", 3.39, "20/10/2023 12:59",
    "https://www.betexplorer.com/football/armenia/premier-league/pyunik-yerevan-ararat-armenia/CxPQkOyq/"
)

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "55").Value = $row55[$i]
    $ws.Range($cols[$i] + "56").Value = $row56[$i]
}

# --- Append two new rows (83 and 84) at the bottom ---

# Copy formatting (style) from the last existing row (82) for the A:E columns
$ws.Range("A82:E82").Copy($ws.Range("A83:E83"))
$ws.Range("A82:E82").Copy($ws.Range("A84:E84"))

# Row 83: Alashkert vs Ararat-Armenia
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "armenia"
$ws.Range("C83").Value = "premier-league"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45260.47916666666
$ws.Range("F83").Value = "Alashkert"
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = "Ararat-Armenia"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 3.38
$ws.Range("K83").Value = "28/11/2023 23:42"
$ws.Range("L83").Value = 3.01
$ws.Range("M83").Value = "30/11/2023 11:21"
$ws.Range("N83").Value = 3.24
$ws.Range("O83").Value = "28/11/2023 23:42"
$ws.Range("P83").Value = 3.37
$ws.Range("Q83").Value = "30/11/2023 11:21"
$ws.Range("R83").Value = 2.05
$ws.Range("S83").Value = "28/11/2023 23:42"
$ws.Range("T83").Value = 2.34
$ws.Range("U83").Value = "30/11/2023 11:21"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/armenia/premier-league/alashkert-ararat-armenia/6uCWgClh/"

# Row 84: BKMA vs Van
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "armenia"
$ws.Range("C84").Value = "premier-league"
$ws.Range("D84").Value = "2023-2024"
$ws.Range("E84").Value = 45260.625
$ws.Range("F84").Value = "BKMA"
$ws.Range("G84").Value = 7
$ws.Range("H84").Value = "Van"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 2.27
$ws.Range("K84").Value = "29/11/2023 03:12"
$ws.Range("L84").Value = 1.67
$ws.Range("M84").Value = "30/11/2023 14:38"
$ws.Range("N84").Value = 3.15
$ws.Range("O84").Value = "29/11/2023 03:12"
$ws.Range("P84").Value = 4.13
$ws.Range("Q84").Value = "30/11/2023 14:38"
$ws.Range("R84").Value = 2.91
$ws.Range("S84").Value = "29/11/2023 03:12"
$ws.Range("T84").Value = 4.63
$ws.Range("U84").Value = "30/11/2023 14:38"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/armenia/premier-league/bkma-van/hdEriAYA/"
